# Update countries & provincias Spain
# This script refreshes the COVID-19 country data table:
#  - Updates totals for Catar and Kuwait (causing them to move up in rank)
#  - Re-sorts the data by "Casos totales" (column B) descending, which re-ranks
#    Catar above Polonia/Rumania and Kuwait above Argelia/Moldavia
#  - Corrects Madagascar's "Casos activos" / "Recuperados" split
#  - Updates the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Catar's row (currently row 36, before re-sort) with refreshed totals ---
$ws.Cells.Item(36, 2).Value = 12564
$ws.Cells.Item(36, 3).Value = 643
$ws.Cells.Item(36, 4).Value = 1243
$ws.Cells.Item(36, 5).Value = 11311
$ws.Cells.Item(36, 6).Value = 72
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 10

# --- Update Kuwait's row (currently row 60, before re-sort) with refreshed totals ---
$ws.Cells.Item(60, 2).Value = 3740
$ws.Cells.Item(60, 3).Value = 300
$ws.Cells.Item(60, 4).Value = 1389
$ws.Cells.Item(60, 5).Value = 2327
$ws.Cells.Item(60, 6).Value = 66
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 24

# --- Correct Madagascar (row 142) active-cases / recovered split ---
$ws.Cells.Item(142, 4).Value = 90
$ws.Cells.Item(142, 5).Value = 38

# --- Re-sort the whole data table (rows 4-216) by "Casos totales" (column B) descending ---
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)

# --- Update the "last updated" timestamp shown in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 13:52"
